$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = -3
    9  = -8
    11 = -1
    17 = 1
    20 = -5
    23 = -3
    24 = 3
    28 = 0
    32 = 0
    36 = -5
    40 = -3
    42 = -2
    46 = -6
    47 = 0
    52 = -7
    56 = 0
    57 = -2
    60 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
